# Updated C3DC phs000466 queries
# Fix the Treatment Agent expression in the TreatmentTab query cell (B5):
# remove the redundant CONCAT(...) wrapper around REPLACE(...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$oldQuery = $treatmentCell.Value2
$newQuery = $oldQuery.Replace(
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))",
    "REPLACE(trt.treatment_agent, ';', ', ')"
)
$treatmentCell.Value2 = $newQuery

# Leave the selection on B2 (matches the saved workbook UI state).
[void]$ws.Range("B2").Select()
